$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final ordered list of employee names for A1:A54 (A1 is the "Nome" header).
# This mirrors the author's edits to the shared-employee list:
#   - "Felipe Rafael Tancredi Pascucci" renamed to "Flavio Henrique Madureira Bergamini"
#   - "Gabriella da Silva Correia" added (after "Gabriel Junior de Souza")
#   - "Guilherme Soares Battaglin" added (after "Guilherme Ribeiro de Melo Yabu")
#   - "Marcelo Bruno Verissimo Mendes Moraes Viegas" removed
# Rows 55:56 stay blank, exactly as before.
$names = @(
  "Nome",
  "Andrea Dalbao",
  "Bruna Eduarda dos Santos Martinez Souza",
  "Bruno Gabriel Nogueira da Silva",
  "Cahue Kokeny Borges Guimarães",
  "Danilo de Jesus Cruz",
  "DIEGO ALVES OPENHEIMER ",
  "Diogenes  Aparecido Rezende",
  "Emanuel Motta Santana Silva",
  "Felipe Gabriel da Cunha",
  "Flavio Henrique Madureira Bergamini",
  "Gabriel Galhato Roriz",
  "Gabriel Junior de Souza",
  "Gabriella da Silva Correia",
  "Giovanni Francez",
  "Guilherme Esquivel dos Santos",
  "Guilherme Ribeiro de Melo Yabu",
  "Guilherme Soares Battaglin",
  "Gustavo de Paiva Caiafa",
  "Gustavo Henrique da Silva Prado",
  "Gustavo Silva Barbosa",
  "Hiago Rangel de Almeida",
  "Hugo Pinheiro Raimundo",
  "João Leonardo Andrade Morganti Silva",
  "João Pedro Giacometti de Souza",
  "João Vitor Alves Lima",
  "Leonardo Everton da Costa",
  "Luana Cristina Cosensa Pierini",
  "Lucas Eduardo Moraes da Silva",
  "Lucas Thalles dos Santos",
  "Matheus Cleber Silva Guerra",
  "Nicole Ribeiro de Paula",
  "Nubia dos Santos Oliveira",
  "Pamela Ferreira Alves Andrelo",
  "Patrick Barnabé Moreira Santos",
  "Patrick Ferreira Araujo",
  "Paulo Gabriel de Freitas Rotundaro",
  "Paulo Sérgio Aquino Ribeiro",
  "Pedro Abritta Reis",
  "Roberto Maia Rodrigues Junior",
  "Roberto Ryan Caldas Ribeiro",
  "Ruan Patrick de Souza",
  "Samuel Alves Brandani Tenório",
  "Samuel da Costa Araujo Nunes",
  "Silas Almeida de Sena",
  "Tarcio Passos Freitas",
  "Tharsis Lamin Lumumba Boa Morte Queiroz",
  "Thiago Agostinho Mem",
  "Valdir Aires Pinheiro Neto",
  "Vanilza Faria de Oliveira",
  "Walber Fellipe de Almeida Rosa",
  "Yasmin Gomes Carpes",
  "Ygor Guilherme Ribeiro Rosa",
  "Vinicius de Castro"
)

# Only rows 11 through 54 actually change value; row 1 ("Nome") and rows 2-10
# are already correct, but writing them again is harmless and keeps this loop simple.
for ($i = 0; $i -lt $names.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $names[$i]
}
